$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header style used by B1:E1 (bold font, centered/top alignment, thin box border)
# by copying E1's formatting into F1 before writing the new header text.
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)

# New header for column F
$ws.Range("F1").Value = "MOP_DEF"

# MOP_DEF values, row 2..9
$defs = @(
    "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']",
    "[]",
    "[]",
    "[]",
    "[]",
    "[]",
    "[]",
    "[]"
)

for ($i = 0; $i -lt $defs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $defs[$i]
}

